$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text assignments (values that Excel will not auto-convert to numbers)
$ws.Range("D2").Value = "41.491.85"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.203.58"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +4.19%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +6.11%  "
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "2.532.37"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("E16").Value = "  +5.47%  "
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "2.190.01"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "41.510.54"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  +20.78%  "
$ws.Range("E26").Value = "  +5.07%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  +6.44%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  +4.39%  "
$ws.Range("E36").Value = "  +7.61%  "
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("E38").Value = "  +6.35%  "
$ws.Range("E39").Value = "  +8.77%  "
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("E41").Value = "  +11.74%  "
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +14.11%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("E51").Value = "  +3.35%  "

# Price cells that look like plain numbers: force text entry via Text format,
# then restore the default style so no extra formatting is left on the cell
# (mirrors how real Excel keeps these columns as text-formatted price strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.36"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.33"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.88"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.68"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0949"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.892"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.64"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.62"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.96"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.88"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.53"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.96"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.80"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0761"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.20"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.62"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.26"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.41"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.93"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.26"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.79"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
